$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record attendance: Lecture on 10/11 (B6) and Team Meeting on 10/13 (D8)
$ws.Range("B6").Value = 1
$ws.Range("D8").Value = 1

# Update the active selection to reflect where the last edit was made
$ws.Range("D8").Select()
